$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.272.04'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.767.87'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.45'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  -1.38%  '
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.52'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("E11").Value = '  +3.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.07'
$ws.Range("E12").Value = '  +3.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0831'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("E14").Value = '  -1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.205.57'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.758.53'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.921'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.210.52'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("E19").Value = '  +3.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.10'
$ws.Range("E20").Value = '  -1.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.11'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.59'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.04'
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.91'
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("E28").Value = '  +13.16%  '
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").Value = '  +1.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.53'
$ws.Range("E31").Value = '  +5.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.81'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.08'
$ws.Range("E33").Value = '  +6.81%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0442'
$ws.Range("E34").Value = '  -4.75%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.51'
$ws.Range("E35").Value = '  +3.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0827'
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.15'
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.90'
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.96'
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.095.78'
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("E50").Value = '  -5.04%  '
$ws.Range("E51").Value = '  +6.82%  '
